# Generate Report for Handback
# Replaces the old handoff/handback UUID-named files with the new ones
# produced by the latest CI run, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldMd1 = "1ace5131-8f07-4735-8afd-ef3048839206.md"
$oldMd2 = "4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.md"
$newMd1 = "beae0205-6e1d-4be2-9975-2e4fbaa5253f.md"
$newMd2 = "ffff813c1e90-df15-475a-8ce6-80aeb00074ce.md"

$newZhCnXlf = "beae0205-6e1d-4be2-9975-2e4fbaa5253f.ccd4440b0fce08797f252197d6afa955b4655a16.zh-cn.xlf"
$newDeDeXlf = "beae0205-6e1d-4be2-9975-2e4fbaa5253f.ccd4440b0fce08797f252197d6afa955b4655a16.de-de.xlf"

$zhCnHandoffTime = "2016-03-12 08:44:50"
$zhCnHandbackTime = "2016-03-12 08:45:07"
$deDeHandoffTime = "2016-03-12 08:44:53"
$deDeHandbackTime = "2016-03-12 08:45:13"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Hyperlinks.Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/1ace5131-8f07-4735-8afd-ef3048839206.md", "", "", $newMd1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.md", "", "", $newMd2)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/1ace5131-8f07-4735-8afd-ef3048839206.md", "", "", $newMd1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/1ace5131-8f07-4735-8afd-ef3048839206.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c28e5164bf5ab32878f89e977d9c8d9a2d623f42/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1ace5131-8f07-4735-8afd-ef3048839206.fe08b98021da4b551880b32dab96ba1d652f2f50.zh-cn.xlf", "", "", $newZhCnXlf)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/337eb03382ba79359ffb49a5e7e245a04475f784/e2e/1ace5131-8f07-4735-8afd-ef3048839206.md", "", "", $newMd1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/29e249a3e20c825db63ad0d8732481be0b307fe3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1ace5131-8f07-4735-8afd-ef3048839206.fe08b98021da4b551880b32dab96ba1d652f2f50.zh-cn.xlf", "", "", $newZhCnXlf)

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.md", "", "", $newMd2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c28e5164bf5ab32878f89e977d9c8d9a2d623f42/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.64c3bc10457d2ae7e5488accef9d1770522e41d0.zh-cn.xlf", "", "", $newZhCnXlf)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/337eb03382ba79359ffb49a5e7e245a04475f784/e2e/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.md", "", "", $newMd2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/29e249a3e20c825db63ad0d8732481be0b307fe3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.64c3bc10457d2ae7e5488accef9d1770522e41d0.zh-cn.xlf", "", "", $newZhCnXlf)

$wsZhCn.Range("E2").Value = $zhCnHandoffTime
$wsZhCn.Range("H2").Value = $zhCnHandbackTime
$wsZhCn.Range("E3").Value = $zhCnHandoffTime
$wsZhCn.Range("H3").Value = $zhCnHandbackTime

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/1ace5131-8f07-4735-8afd-ef3048839206.md", "", "", $newMd1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/1ace5131-8f07-4735-8afd-ef3048839206.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/054b6622107f28cbb7ae1d38512c3292ed8b31cb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1ace5131-8f07-4735-8afd-ef3048839206.fe08b98021da4b551880b32dab96ba1d652f2f50.de-de.xlf", "", "", $newDeDeXlf)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/3866860413470ee970a1e4894c3ce2ccf91b27d7/e2e/1ace5131-8f07-4735-8afd-ef3048839206.md", "", "", $newMd1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9d834b6f3218155ad8be71a4f558dc6c5da22c58/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1ace5131-8f07-4735-8afd-ef3048839206.fe08b98021da4b551880b32dab96ba1d652f2f50.de-de.xlf", "", "", $newDeDeXlf)

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.md", "", "", $newMd2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/696c95751374da0531a9c1e131c9071384cc71dd/e2e/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/054b6622107f28cbb7ae1d38512c3292ed8b31cb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.64c3bc10457d2ae7e5488accef9d1770522e41d0.de-de.xlf", "", "", $newDeDeXlf)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/3866860413470ee970a1e4894c3ce2ccf91b27d7/e2e/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.md", "", "", $newMd2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9d834b6f3218155ad8be71a4f558dc6c5da22c58/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.64c3bc10457d2ae7e5488accef9d1770522e41d0.de-de.xlf", "", "", $newDeDeXlf)

$wsDeDe.Range("E2").Value = $deDeHandoffTime
$wsDeDe.Range("H2").Value = $deDeHandbackTime
$wsDeDe.Range("E3").Value = $deDeHandoffTime
$wsDeDe.Range("H3").Value = $deDeHandbackTime
